$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 148, shifting existing rows
# (old 148-171) down to (149-172).
$ws.Rows("148").Insert()

# Populate the newly inserted row 148 with the new weekly data point.
$ws.Range("A148").Value = 3
$ws.Range("B148").Value = "Femacal de La Calera"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44776
$ws.Range("E148").Value = 5
$ws.Range("F148").Value = 100112026
$ws.Range("G148").Value = "Haba"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 85
$ws.Range("K148").Value = 16000
$ws.Range("L148").Value = 17000
$ws.Range("M148").Value = 16529
$ws.Range("N148").Value = "$/saco 25 kilos"
$ws.Range("O148").Value = "Provincia de Limarí"
$ws.Range("P148").Value = 661
$ws.Range("Q148").Value = 25
$ws.Range("R148").Value = "Hortaliza"
